$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the pipeline counts in column C; dependent ratio formulas in column D
# (C3/C2, C4/C3, C5/C4, C6/C5, C7/C6, C8/C7) recalc automatically.
$ws.Range("C2").Value = 17686161
$ws.Range("C3").Value = 5562533
$ws.Range("C4").Value = 918988
$ws.Range("C5").Value = 756588
$ws.Range("C6").Value = 44735
$ws.Range("C7").Value = 38989
$ws.Range("C8").Value = 433

# Move the selection to C6 (single cell) like in the saved view state.
$null = $ws.Range("C6").Select()
